# REPORTGEN-1102: part 1, added and removed counts missing when no previous
# snapshot selected.
#
# The RepGen "evolution" table placeholders (QUALITY_STANDARDS_EVOLUTION and
# RULES_LIST_STATISTICS_RATIO) need an extra EVOLUTION=true flag appended to
# their instruction text so that the report generator includes the
# added/removed counts even when no previous snapshot was selected.

$wb = $excel.ActiveWorkbook

# Summary sheet: the "Findings summary" evolution table.
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B14").Value = "RepGen:TABLE;QUALITY_STANDARDS_EVOLUTION;STD=AIP-MIPS-Reduction,HEADER=NO,EVOLUTION=true"

# MIPS-ALGORITHMIC-COST sheet: rules statistics ratio evolution table.
$wsAlgo = $wb.Worksheets.Item("MIPS-ALGORITHMIC-COST")
$wsAlgo.Range("A3").Value = "RepGen:TABLE;RULES_LIST_STATISTICS_RATIO;METRICS=AIP-MIPS-ALGORITHMIC-COST,DESC=true,HEADER=NO,EVOLUTION=true"

# MIPS-DATA-ACCESS sheet: rules statistics ratio evolution table.
$wsData = $wb.Worksheets.Item("MIPS-DATA-ACCESS")
$wsData.Range("A3").Value = "RepGen:TABLE;RULES_LIST_STATISTICS_RATIO;METRICS=AIP-MIPS-DATA-ACCESS,DESC=true,HEADER=NO,EVOLUTION=true"

# MIPS-TRANS-FAIL sheet: rules statistics ratio evolution table.
$wsTrans = $wb.Worksheets.Item("MIPS-TRANS-FAIL")
$wsTrans.Range("A3").Value = "RepGen:TABLE;RULES_LIST_STATISTICS_RATIO;METRICS=AIP-MIPS-TRANSACTIONS-FAILURE,DESC=true,HEADER=NO,EVOLUTION=true"
